# Updates numeric Leve price/profit columns (H,I,J,K,L,M,N) across multiple
# sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR) to reflect refreshed market data.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 43
$ws.Range("I11").Value = 43
$ws.Range("K11").Value = 43
$ws.Range("M11").Value = 97
$ws.Range("H15").Value = 890.2727
$ws.Range("I15").Value = 890.2727
$ws.Range("K15").Value = 2670.8181
$ws.Range("M15").Value = -2501.8181
$ws.Range("H18").Value = 981.4286
$ws.Range("I18").Value = 811.6667
$ws.Range("K18").Value = 811.6667
$ws.Range("M18").Value = -527.6667
$ws.Range("H38").Value = 3235.9167
$ws.Range("I38").Value = 1687.5714
$ws.Range("J38").Value = 5403.6
$ws.Range("K38").Value = 5062.7142
$ws.Range("L38").Value = 16210.8
$ws.Range("M38").Value = -4690.7142
$ws.Range("N38").Value = -16954.8
$ws.Range("H53").Value = 215.66667
$ws.Range("I53").Value = 117.625
$ws.Range("J53").Value = 1000
$ws.Range("K53").Value = 117.625
$ws.Range("L53").Value = 1000
$ws.Range("M53").Value = 519.375
$ws.Range("N53").Value = -2274
$ws.Range("H58").Value = 2233
$ws.Range("J58").Value = 2677.0908
$ws.Range("L58").Value = 8031.2724
$ws.Range("N58").Value = -8331.2724
$ws.Range("H62").Value = 6300.5
$ws.Range("I62").Value = 3661.6667
$ws.Range("K62").Value = 3661.6667
$ws.Range("M62").Value = -3037.6667
$ws.Range("H65").Value = 6300.5
$ws.Range("I65").Value = 3661.6667
$ws.Range("K65").Value = 18308.3335
$ws.Range("M65").Value = -15188.3335
$ws.Range("H70").Value = 10000
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30540
$ws.Range("H73").Value = 10000
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -31872
$ws.Range("H75").Value = 40167.332
$ws.Range("J75").Value = 35001
$ws.Range("L75").Value = 35001
$ws.Range("N75").Value = -36873
$ws.Range("H76").Value = 3981.3333
$ws.Range("I76").Value = 3944
$ws.Range("K76").Value = 3944
$ws.Range("M76").Value = -3629
$ws.Range("H78").Value = 40167.332
$ws.Range("J78").Value = 35001
$ws.Range("L78").Value = 105003
$ws.Range("N78").Value = -114363
$ws.Range("H79").Value = 3981.3333
$ws.Range("I79").Value = 3944
$ws.Range("K79").Value = 3944
$ws.Range("M79").Value = -2852
$ws.Range("H88").Value = 759.375
$ws.Range("I88").Value = 817
$ws.Range("J88").Value = 701.75
$ws.Range("K88").Value = 817
$ws.Range("L88").Value = 701.75
$ws.Range("M88").Value = -411
$ws.Range("N88").Value = -1513.75
$ws.Range("H91").Value = 759.375
$ws.Range("I91").Value = 817
$ws.Range("J91").Value = 701.75
$ws.Range("K91").Value = 817
$ws.Range("L91").Value = 701.75
$ws.Range("M91").Value = 587
$ws.Range("N91").Value = -3509.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3576001.5
$ws.Range("I32").Value = 5255.9375
$ws.Range("K32").Value = 5255.9375
$ws.Range("M32").Value = -4968.9375
$ws.Range("H45").Value = 2597.2917
$ws.Range("I45").Value = 2067.4
$ws.Range("J45").Value = 3480.4443
$ws.Range("K45").Value = 2067.4
$ws.Range("L45").Value = 3480.4443
$ws.Range("M45").Value = -1690.4
$ws.Range("N45").Value = -4234.4443
$ws.Range("H110").Value = 531.1
$ws.Range("I110").Value = 516.375
$ws.Range("K110").Value = 516.375
$ws.Range("M110").Value = 1528.625
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7048.625
$ws.Range("J107").Value = 8100
$ws.Range("L107").Value = 8100
$ws.Range("N107").Value = -11940

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6341.3335
$ws.Range("I31").Value = 995
$ws.Range("K31").Value = 995
$ws.Range("M31").Value = -700
$ws.Range("H34").Value = 6341.3335
$ws.Range("I34").Value = 995
$ws.Range("K34").Value = 995
$ws.Range("M34").Value = -793
$ws.Range("H58").Value = 3820.2
$ws.Range("I58").Value = 3315.5715
$ws.Range("K58").Value = 3315.5715
$ws.Range("M58").Value = -3112.5715
$ws.Range("H134").Value = 2000
$ws.Range("I134").Value = 3000
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 9000
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -6465
$ws.Range("N134").Value = -8070
$ws.Range("H136").Value = 3820.2
$ws.Range("I136").Value = 3315.5715
$ws.Range("K136").Value = 9946.7145
$ws.Range("M136").Value = -7396.7145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 735.875
$ws.Range("J5").Value = 778.44446
$ws.Range("L5").Value = 2335.33338
$ws.Range("N5").Value = -2559.33338
$ws.Range("H44").Value = 499.75
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 499.75
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 1499.25
$ws.Range("M44").Value = ""
$ws.Range("N44").Value = -2295.25
$ws.Range("H50").Value = 568.3333
$ws.Range("I50").Value = 352.5
$ws.Range("K50").Value = 1057.5
$ws.Range("M50").Value = -576.5
$ws.Range("H53").Value = 568.3333
$ws.Range("I53").Value = 352.5
$ws.Range("K53").Value = 1057.5
$ws.Range("M53").Value = -576.5
$ws.Range("H68").Value = 1996
$ws.Range("I68").Value = 1625
$ws.Range("J68").Value = 2243.3333
$ws.Range("K68").Value = 4875
$ws.Range("L68").Value = 6729.999899999999
$ws.Range("M68").Value = -4064
$ws.Range("N68").Value = -8351.999899999999
$ws.Range("H71").Value = 1996
$ws.Range("I71").Value = 1625
$ws.Range("J71").Value = 2243.3333
$ws.Range("K71").Value = 14625
$ws.Range("L71").Value = 20189.9997
$ws.Range("M71").Value = -10569
$ws.Range("N71").Value = -28301.9997
$ws.Range("H121").Value = 457.6
$ws.Range("J121").Value = 983
$ws.Range("L121").Value = 2949
$ws.Range("N121").Value = -5569
$ws.Range("H131").Value = 1620.409
$ws.Range("I131").Value = 913.3
$ws.Range("J131").Value = 2209.6667
$ws.Range("K131").Value = 2739.9
$ws.Range("L131").Value = 6629.000100000001
$ws.Range("M131").Value = 2300.1
$ws.Range("N131").Value = -16709.0001
$ws.Range("H135").Value = 735.875
$ws.Range("J135").Value = 778.44446
$ws.Range("L135").Value = 7006.00014
$ws.Range("N135").Value = -12076.00014
$ws.Range("H138").Value = 4088.4546
$ws.Range("I138").Value = 2998
$ws.Range("K138").Value = 8994
$ws.Range("M138").Value = -3854
$ws.Range("H139").Value = 3738
$ws.Range("I139").Value = 3206
$ws.Range("K139").Value = 9618
$ws.Range("M139").Value = -4478

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 65000
$ws.Range("J69").Value = 65000
$ws.Range("L69").Value = 65000
$ws.Range("N69").Value = -66622
$ws.Range("H72").Value = 65000
$ws.Range("J72").Value = 65000
$ws.Range("L72").Value = 195000
$ws.Range("N72").Value = -203112
$ws.Range("H93").Value = 986.625
$ws.Range("I93").Value = 838.6
$ws.Range("K93").Value = 838.6
$ws.Range("M93").Value = 409.4
$ws.Range("H100").Value = 6808.231
$ws.Range("I100").Value = 4800.6
$ws.Range("J100").Value = 8063
$ws.Range("K100").Value = 4800.6
$ws.Range("L100").Value = 8063
$ws.Range("M100").Value = -4259.6
$ws.Range("N100").Value = -9145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 51999.5
$ws.Range("J54").Value = 99999
$ws.Range("L54").Value = 99999
$ws.Range("N54").Value = -101039
$ws.Range("H62").Value = 12333.333
$ws.Range("I62").Value = 9500
$ws.Range("J62").Value = 12900
$ws.Range("K62").Value = 9500
$ws.Range("L62").Value = 12900
$ws.Range("M62").Value = -8876
$ws.Range("N62").Value = -14148
$ws.Range("H65").Value = 12333.333
$ws.Range("I65").Value = 9500
$ws.Range("J65").Value = 12900
$ws.Range("K65").Value = 47500
$ws.Range("L65").Value = 64500
$ws.Range("M65").Value = -44380
$ws.Range("N65").Value = -70740
$ws.Range("H80").Value = 87454
$ws.Range("J80").Value = 87454
$ws.Range("L80").Value = 87454
$ws.Range("N80").Value = -89450
$ws.Range("H83").Value = 87454
$ws.Range("J83").Value = 87454
$ws.Range("L83").Value = 262362
$ws.Range("N83").Value = -272346
$ws.Range("H132").Value = 2247
$ws.Range("I132").Value = 2277.875
$ws.Range("K132").Value = 6833.625
$ws.Range("M132").Value = -4303.625
$ws.Range("H136").Value = 4254.95
$ws.Range("I136").Value = 3299.2
$ws.Range("K136").Value = 9897.599999999999
$ws.Range("M136").Value = -7347.599999999999

